$d = $word.ActiveDocument

function Replace-Text($old, $new) {
  $range = $d.Content
  $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
  if (-not $found) {
    Write-Host "NOT FOUND: $old"
  }
}

Replace-Text "2024-06-24 Monday" "2024-06-25 Tuesday"
Replace-Text "31-27=4" "67+8=75"
Replace-Text "18+60=78" "65+15=80"
Replace-Text "92+6=98" "55-6=49"
Replace-Text "58+21=79" "73-67=6"
Replace-Text "17+33=50" "35+25=60"
Replace-Text "85-33=52" "84-70=14"
Replace-Text "28+11=39" "82-2=80"
Replace-Text "7+24=31" "38-15=23"
Replace-Text "83+13=96" "61+5=66"
Replace-Text "79+2=81" "42+25=67"
Replace-Text "94-67=27" "17+76=93"
Replace-Text "72-11=61" "31+18=49"
Replace-Text "79-60=19" "4+88=92"
Replace-Text "81-10=71" "55-24=31"
Replace-Text "15+61=76" "32-29=3"
Replace-Text "25-14=11" "43+52=95"
Replace-Text "8+79=87" "25+44=69"
Replace-Text "91-84=7" "80-71=9"
Replace-Text "1+92=93" "34+44=78"
Replace-Text "89-3=86" "8+6=14"
Replace-Text "27+63=90" "48-44=4"
Replace-Text "66-29=37" "55-37=18"
Replace-Text "77-20=57" "31+49=80"
Replace-Text "26+37=63" "11+83=94"
Replace-Text "35+6=41" "34+20=54"
Replace-Text "21+14=35" "6+53=59"
Replace-Text "18+57=75" "7+34=41"
Replace-Text "55-3=52" "94-19=75"
Replace-Text "89-37=52" "78+17=95"
Replace-Text "27+69=96" "20-6=14"
Replace-Text "92-64=28" "90-74=16"
Replace-Text "85-61=24" "84+2=86"
Replace-Text "69-19=50" "33-21=12"
Replace-Text "4+78=82" "31+54=85"
Replace-Text "6+69=75" "68+11=79"
Replace-Text "67-65=2" "9+67=76"
Replace-Text "56+14=70" "84-50=34"
Replace-Text "62-31=31" "45+28=73"
Replace-Text "32+36=68" "41+37=78"
Replace-Text "93-57=36" "83-41=42"
Replace-Text "62+19=81" "71-52=19"
Replace-Text "66-40=26" "34+24=58"
Replace-Text "10+26=36" "28+46=74"
Replace-Text "93-50=43" "3+72=75"
Replace-Text "22+72=94" "14+66=80"
Replace-Text "72+18=90" "9+35=44"
Replace-Text "74-3=71" "91-50=41"
Replace-Text "36+46=82" "26+9=35"
Replace-Text "1+16=17" "87-77=10"
Replace-Text "6+56=62" "28+28=56"
Replace-Text "34+56=90" "73+14=87"
Replace-Text "34+10=44" "32+50=82"
Replace-Text "67+0=67" "26+20=46"
Replace-Text "35+16=51" "46+11=57"
Replace-Text "43+17=60" "30-4=26"
Replace-Text "26+14=40" "23+20=43"
Replace-Text "65+25=90" "51-38=13"
Replace-Text "26+36=62" "96-71=25"
Replace-Text "56+40=96" "63+29=92"
Replace-Text "71-32=39" "32+66=98"
Replace-Text "32-7=25" "12-7=5"
Replace-Text "33+16=49" "82-58=24"
Replace-Text "40+4=44" "25+72=97"
Replace-Text "84-44=40" "92-19=73"
Replace-Text "40-32=8" "31+37=68"
Replace-Text "63-26=37" "88-48=40"
Replace-Text "98-74=24" "33+24=57"
Replace-Text "32-26=6" "68-29=39"
Replace-Text "10+6=16" "98-13=85"
Replace-Text "57+40=97" "84-81=3"
Replace-Text "16+81=97" "32+6=38"
Replace-Text "15+33=48" "36+24=60"
Replace-Text "66+26=92" "53+26=79"
Replace-Text "42+22=64" "11+34=45"
Replace-Text "58-31=27" "25-13=12"
Replace-Text "86-78=8" "65-10=55"
Replace-Text "29+35=64" "31+46=77"
Replace-Text "46+51=97" "97-63=34"
Replace-Text "58+7=65" "16+50=66"
Replace-Text "11+87=98" "83-40=43"
Replace-Text "40+23=63" "19+58=77"
Replace-Text "96-12=84" "83+1=84"
Replace-Text "63-55=8" "5+13=18"
Replace-Text "68-48=20" "71-43=28"
Replace-Text "90-63=27" "95-43=52"
Replace-Text "69-53=16" "52-30=22"
Replace-Text "51+29=80" "11+53=64"
Replace-Text "7+89=96" "30+20=50"
Replace-Text "73+9=82" "24+63=87"
Replace-Text "20+56=76" "13+35=48"
Replace-Text "86-14=72" "60+15=75"
Replace-Text "87-49=38" "94-38=56"
Replace-Text "10-7=3" "4+28=32"
Replace-Text "98-49=49" "57+26=83"
Replace-Text "34-29=5" "14+15=29"
Replace-Text "66+18=84" "24+24=48"
Replace-Text "31-2=29" "98-88=10"
Replace-Text "47+33=80" "95-33=62"
Replace-Text "75+6=81" "44+55=99"
Replace-Text "27+9=36" "1+19=20"
